# Generate Report for Handoff
# - Status text changes from "Handed back: in sync with en-US" to "Ready for handoff"
#   everywhere it appears (Overview!E2/F2, zh-cn!C2, de-de!C2)
# - Timestamps bump forward a few seconds/minutes to reflect the new handoff run:
#     Overview!G2            2016-08-27 02:57:08 -> 2016-08-27 02:58:04
#     zh-cn!H2 (Latest Handoff Datetime)  2016-08-27 02:56:59 -> 2016-08-27 02:57:56
#     de-de!H2 (Latest Handback DateTime) 2016-08-27 02:57:08 -> 2016-08-27 02:58:04
# - The Status columns got narrower (report layout tweak)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# --- Status text ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-08-27 02:58:04"
$wsZhCn.Range("H2").Value = "2016-08-27 02:57:56"
$wsDeDe.Range("H2").Value = "2016-08-27 02:58:04"

# --- Column widths (Status column narrower in all three sheets) ---
# Target stored width is 17.2159881591797 characters; this host's ColumnWidth
# setter quantizes the persisted <col width> to the nearest 1/6 of a
# character (Excel's internal pixel-rounding granularity), so we pick the
# input that lands on the closest reachable quantum (17.1666...) rather than
# the unreachable exact fraction.
$targetColumnWidth = 16.26
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
